{"js": "// Append, after the last paragraph in the document body\n// (\"Reiterate how it is on track.\"), a blank paragraph followed by a\n// new (non-list, \"Normal\" style) paragraph containing:\n//   \"Process generally completed dimensions \"\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Blank separator paragraph, matching the existing blank-line style\n// already used elsewhere in the document (e.g. right after the title).\nconst blankParagraph = lastParagraph.insertParagraph(\"\", \"After\");\nblankParagraph.style = \"Normal\";\nawait context.sync();\n\n// New paragraph with the actual text. Leave its style untouched after\n// creation (it already inherits \"Normal\" from blankParagraph), so no\n// redundant pStyle gets written out.\nblankParagraph.insertParagraph(\"Process generally completed dimensions \", \"After\");\n\nawait context.sync();\n", "ps1": "# Append, after the last paragraph in the document\n# (\"Reiterate how it is on track.\"), a blank paragraph followed by a\n# new (non-list, \"Normal\" style) paragraph containing:\n#   \"Process generally completed dimensions \"\n\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n$endOfDoc = $lastParagraph.Range\n$endOfDoc.Collapse(0)          # wdCollapseEnd\n$endOfDoc.InsertParagraphAfter()\n\n# The paragraph we just created is now the last one in the document.\n$blankParagraph = $d.Paragraphs.Last\n$blankParagraph.Style = \"Normal\"\n\n$blankRange = $blankParagraph.Range\n$blankRange.Collapse(0)\n$blankRange.InsertParagraphAfter()\n\n# New paragraph holding the actual text.\n$newParagraph = $d.Paragraphs.Last\n$newParagraph.Range.Text = \"Process generally completed dimensions \"\n"}
